$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (POLID 123456) currently has PAYMENT STATUS "E" (Error/Excluded).
# Update it to "M" to include it in the comparison as an error row.
$ws.Range("D2").Value = "M"
